# Add an "ID" column (column B) next to the existing "Clients" column:
# a styled bold header row (yellow / green fill) and a centered,
# text-formatted, green-filled numeric id for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# B1 gets the new "ID" label
$ws.Range("B1").Value = "ID"

# A1 "Clients" header style: bold 12pt font, solid yellow fill, centered
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Font.Size = 12
$a1.Interior.Color = 65535
$a1.Interior.PatternColor = 65535
$a1.HorizontalAlignment = -4108

# B1 "ID" header style: bold 12pt font, solid green fill, centered
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Font.Size = 12
$b1.Interior.Color = 3385600
$b1.Interior.PatternColor = 32768
$b1.HorizontalAlignment = -4108

$ws.Rows.Item(1).RowHeight = 15

# ---- Data rows (2-7) ----
# New numeric "ID" values 1-6 in column B
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 6

# Column A data style: yellow fill (font/format unchanged) - build it once
# on A2 then stamp the rest of the column via copy/paste of formats only,
# so every A2:A7 cell shares a single style record.
$a2 = $ws.Range("A2")
$a2.Interior.Color = 65535
$a2.Interior.PatternColor = 65535
$a2.Copy() | Out-Null
$ws.Range("A3:A7").PasteSpecial(-4122) | Out-Null

# Column B data style: text number format, green fill, centered - build it
# once on B2 then stamp the rest of the column the same way.
$b2 = $ws.Range("B2")
$b2.NumberFormat = "@"
$b2.Interior.Color = 3385600
$b2.Interior.PatternColor = 32768
$b2.HorizontalAlignment = -4108
$b2.VerticalAlignment = -4108
$b2.Copy() | Out-Null
$ws.Range("B3:B7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---- Final selection, mirroring the "run and view" command context ----
$ws.Range("I7").Select() | Out-Null
